# Update the NATMI ligand-receptor pair table (Lrfn3-Lrfn3) with refreshed
# TPM-derived values. The underlying "Ligand/Receptor average expression
# value" for the ECs and MuSCs clusters changed, which cascades into every
# dependent column (totals, specificities, and edge weights) across all
# nine sending/target cluster combinations (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.167796
$ws.Range("H2").Value = 0.5033880000000001
$ws.Range("I2").Value = 0.03973966974868284
$ws.Range("J2").Value = 0.03973966974868284
$ws.Range("M2").Value = 0.167796
$ws.Range("N2").Value = 0.5033880000000001
$ws.Range("O2").Value = 0.03973966974868284
$ws.Range("P2").Value = 0.03973966974868284
$ws.Range("Q2").Value = 0.02815549761600001
$ws.Range("R2").Value = 0.2533994785440001
$ws.Range("S2").Value = 0.001579241351734378
$ws.Range("T2").Value = 0.001579241351734378
$ws.Range("G3").Value = 0.167796
$ws.Range("H3").Value = 0.5033880000000001
$ws.Range("I3").Value = 0.03973966974868284
$ws.Range("J3").Value = 0.03973966974868284
$ws.Range("O3").Value = 0.751323601750387
$ws.Range("P3").Value = 0.751323601750387
$ws.Range("Q3").Value = 0.5323116677040001
$ws.Range("R3").Value = 4.790805009336001
$ws.Range("S3").Value = 0.02985735180795129
$ws.Range("T3").Value = 0.02985735180795129
$ws.Range("G4").Value = 0.167796
$ws.Range("H4").Value = 0.5033880000000001
$ws.Range("I4").Value = 0.03973966974868284
$ws.Range("J4").Value = 0.03973966974868284
$ws.Range("M4").Value = 0.8822103333333334
$ws.Range("N4").Value = 2.646631
$ws.Range("O4").Value = 0.2089367285009301
$ws.Range("P4").Value = 0.2089367285009301
$ws.Range("Q4").Value = 0.148031365092
$ws.Range("R4").Value = 1.332282285828
$ws.Range("S4").Value = 0.00830307658899717
$ws.Range("T4").Value = 0.00830307658899717
$ws.Range("I5").Value = 0.751323601750387
$ws.Range("J5").Value = 0.751323601750387
$ws.Range("M5").Value = 0.167796
$ws.Range("N5").Value = 0.5033880000000001
$ws.Range("O5").Value = 0.03973966974868284
$ws.Range("P5").Value = 0.03973966974868284
$ws.Range("Q5").Value = 0.5323116677040001
$ws.Range("R5").Value = 4.790805009336001
$ws.Range("S5").Value = 0.02985735180795129
$ws.Range("T5").Value = 0.02985735180795129
$ws.Range("I6").Value = 0.751323601750387
$ws.Range("J6").Value = 0.751323601750387
$ws.Range("O6").Value = 0.751323601750387
$ws.Range("P6").Value = 0.751323601750387
$ws.Range("S6").Value = 0.5644871545471741
$ws.Range("T6").Value = 0.5644871545471741
$ws.Range("I7").Value = 0.751323601750387
$ws.Range("J7").Value = 0.751323601750387
$ws.Range("M7").Value = 0.8822103333333334
$ws.Range("N7").Value = 2.646631
$ws.Range("O7").Value = 0.2089367285009301
$ws.Range("P7").Value = 0.2089367285009301
$ws.Range("Q7").Value = 2.798701123998
$ws.Range("R7").Value = 25.188310115982
$ws.Range("S7").Value = 0.1569790953952615
$ws.Range("T7").Value = 0.1569790953952615
$ws.Range("G8").Value = 0.8822103333333334
$ws.Range("H8").Value = 2.646631
$ws.Range("I8").Value = 0.2089367285009301
$ws.Range("J8").Value = 0.2089367285009301
$ws.Range("M8").Value = 0.167796
$ws.Range("N8").Value = 0.5033880000000001
$ws.Range("O8").Value = 0.03973966974868284
$ws.Range("P8").Value = 0.03973966974868284
$ws.Range("Q8").Value = 0.148031365092
$ws.Range("R8").Value = 1.332282285828
$ws.Range("S8").Value = 0.00830307658899717
$ws.Range("T8").Value = 0.00830307658899717
$ws.Range("G9").Value = 0.8822103333333334
$ws.Range("H9").Value = 2.646631
$ws.Range("I9").Value = 0.2089367285009301
$ws.Range("J9").Value = 0.2089367285009301
$ws.Range("O9").Value = 0.751323601750387
$ws.Range("P9").Value = 0.751323601750387
$ws.Range("Q9").Value = 2.798701123998
$ws.Range("R9").Value = 25.188310115982
$ws.Range("S9").Value = 0.1569790953952615
$ws.Range("T9").Value = 0.1569790953952615
$ws.Range("G10").Value = 0.8822103333333334
$ws.Range("H10").Value = 2.646631
$ws.Range("I10").Value = 0.2089367285009301
$ws.Range("J10").Value = 0.2089367285009301
$ws.Range("M10").Value = 0.8822103333333334
$ws.Range("N10").Value = 2.646631
$ws.Range("O10").Value = 0.2089367285009301
$ws.Range("P10").Value = 0.2089367285009301
$ws.Range("Q10").Value = 0.7782950722401113
$ws.Range("R10").Value = 7.004655650161001
$ws.Range("S10").Value = 0.04365455651667137
$ws.Range("T10").Value = 0.04365455651667137
